$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 02:59:49"
$wsZhCn.Range("H2").Value = "2016-03-13 03:00:14"

# de-de sheet: update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 02:59:54"
$wsDeDe.Range("H2").Value = "2016-03-13 03:00:20"
